$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("13-10-2021", "14-10-2021", "15-10-2021", "18-10-2021")

$row = 197
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = 3.25
    $row++
}
